$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Russian localization column (C) used "{...FactionDef.Demonym}" in a number
# of translated strings where the English source had already been updated to
# use "{...FactionDef.Name}". Bring the Russian column in line by replacing
# every remaining "Demonym" token with "Name" in column C only.
$col = $ws.Range("C:C")
$col.Select()
$col.Replace("Demonym", "Name") | Out-Null
